$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the value to be stored as literal text (apostrophe-prefix),
    # avoiding Excel's automatic text-to-number inference, then strip
    # the resulting "quote prefix" style so formatting matches the
    # original (unstyled) cells.
    $rng = $ws.Range($addr)
    $rng.Formula = "'" + $val
    $rng.Style = "Normal"
}

$sub6 = [char]0x2086
$row46Price = "0.0$($sub6)0108"

# Row 2 - Bitcoin
Set-TextValue "D2" "26.213.58"
Set-TextValue "E2" "  +1.54%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.644.40"
Set-TextValue "E3" "  +0.27%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  -0.08%  "

# Row 5 - BNB
Set-TextValue "D5" "216.83"
Set-TextValue "E5" "  +0.41%  "

# Row 6 - XRP
Set-TextValue "D6" "0.505"
Set-TextValue "E6" "  +0.60%  "

# Row 7 - USDC
Set-TextValue "E7" "  -0.11%  "

# Row 8 - Cardano
Set-TextValue "E8" "  +0.12%  "

# Row 9 - Dogecoin
Set-TextValue "E9" "  +0.15%  "

# Row 10 - Solana
Set-TextValue "D10" "19.92"
Set-TextValue "E10" "  +1.13%  "

# Row 11 - TRON
Set-TextValue "E11" "  +0.27%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "1.871.88"
Set-TextValue "E12" "  +0.33%  "

# Row 13 - Polkadot
Set-TextValue "E13" "  +0.71%  "

# Row 14 - WrappedEther
Set-TextValue "D14" "1.634.84"
Set-TextValue "E14" "  -0.30%  "

# Row 15 - Polygon
Set-TextValue "E15" "  -2.73%  "

# Row 16 - ShibaInu
Set-TextValue "E16" "  -0.21%  "

# Row 17 - Litecoin
Set-TextValue "E17" "  +0.21%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "26.214.63"
Set-TextValue "E18" "  +1.43%  "

# Row 19 - Dai
Set-TextValue "E19" "  -0.11%  "

# Row 20 - Uniswap
Set-TextValue "E20" "  -0.90%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "195.12"
Set-TextValue "E21" "  +1.35%  "

# Row 22 - Avalanche
Set-TextValue "D22" "10.05"
Set-TextValue "E22" "  +0.74%  "

# Row 23 - Chainlink
Set-TextValue "E23" "  -0.39%  "

# Row 24 - Toncoin
Set-TextValue "E24" "  -3.00%  "

# Row 25 - BinanceUSD -> Monero (rows 25 and 26 swap coins)
Set-TextValue "B25" "Monero"
Set-TextValue "C25" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D25" "143.03"
Set-TextValue "E25" "  +0.52%  "

# Row 26 - Monero -> BinanceUSD
Set-TextValue "B26" "BinanceUSD"
Set-TextValue "C26" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D26" "1.00"
Set-TextValue "E26" "  -0.17%  "

# Row 27 - Stellar
Set-TextValue "E27" "  +0.95%  "

# Row 28 - Cosmos
Set-TextValue "E28" "  +0.14%  "

# Row 29 - EthereumClassic
Set-TextValue "D29" "15.63"
Set-TextValue "E29" "  +0.75%  "

# Row 30 - PancakeSwap
Set-TextValue "E30" "  +0.64%  "

# Row 31 - Hedera
Set-TextValue "E31" "  +2.02%  "

# Row 32 - InternetComputer(DFINITY)
Set-TextValue "E32" "  +0.33%  "

# Row 33 - Filecoin
Set-TextValue "E33" "  +0.33%  "

# Row 34 - LidoDAOToken
Set-TextValue "E34" "  +1.82%  "

# Row 35 - HuobiToken
Set-TextValue "D35" "2.41"
Set-TextValue "E35" "  +1.29%  "

# Row 36 - ARBITRUM
Set-TextValue "E36" "  +0.41%  "

# Row 37 - Maker
Set-TextValue "D37" "1.136.29"

# Row 38 - ImmutableX
Set-TextValue "E38" "  +1.50%  "

# Row 40 - VeChain
Set-TextValue "E40" "  +1.06%  "

# Row 41 - PaxDollar
Set-TextValue "E41" "  -0.17%  "

# Row 42 - Quant
Set-TextValue "D42" "100.30"
Set-TextValue "E42" "  -0.44%  "

# Row 43 - FraxShare
Set-TextValue "D43" "5.51"
Set-TextValue "E43" "  -1.28%  "

# Row 44 - TrustWalletToken
Set-TextValue "D44" "0.800"
Set-TextValue "E44" "  -0.64%  "

# Row 45 - RocketPoolETH
Set-TextValue "D45" "1.781.40"
Set-TextValue "E45" "  +0.36%  "

# Row 46 - BabyDogeCoin
Set-TextValue "D46" $row46Price
Set-TextValue "E46" "  -3.25%  "

# Row 47 - Aave
Set-TextValue "D47" "56.65"
Set-TextValue "E47" "  +2.31%  "

# Row 48 - RenderToken
Set-TextValue "E48" "  +3.59%  "

# Row 49 - Cronos
Set-TextValue "D49" "0.0517"
Set-TextValue "E49" "  +3.00%  "

# Row 50 - Mantle
Set-TextValue "E50" "  +0.25%  "

# Row 51 - EnergySwap
Set-TextValue "E51" "  +2.71%  "
